# "update core tms 1&2"
# Sheet1!A2 ("id" column) holds a generated test-run identifier
# (e.g. CA-5XAUJ6J3). Re-running the id generator produced a new
# value, CA-TVP8RWH8, which replaces the previous one in A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "CA-TVP8RWH8"
